$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 169. This shifts the existing rows 169-179
# down to 170-180 and extends the used range to row 180.
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with the new weekly record.
$ws.Cells.Item(169, 1).Value = 10
$ws.Cells.Item(169, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(169, 3).Value = "La Araucanía"
$ws.Cells.Item(169, 4).Value = 44610
$ws.Cells.Item(169, 5).Value = 9
$ws.Cells.Item(169, 6).Value = 100112005
$ws.Cells.Item(169, 7).Value = "Puerro"
$ws.Cells.Item(169, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 30
$ws.Cells.Item(169, 11).Value = 10000
$ws.Cells.Item(169, 12).Value = 10000
$ws.Cells.Item(169, 13).Value = 10000
$ws.Cells.Item(169, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(169, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(169, 16).Value = 833
$ws.Cells.Item(169, 17).Value = 12
$ws.Cells.Item(169, 18).Value = "Hortaliza"
